$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Chequia (row 45) figures ---
$ws.Range("B45").Value = 7408
$ws.Range("C45").Value = 4
$ws.Range("D45").Value = 2600
$ws.Range("E45").Value = 4587

# --- Update Moldavia (row 58) figures ---
$ws.Range("D58").Value = 925
$ws.Range("E58").Value = 2382
$ws.Range("G58").Value = 5
$ws.Range("H58").Value = 101

# --- Etiopia moves above Madagascar with refreshed figures (rows 141-142) ---
$ws.Rows("141").Insert()

$ws.Range("A141").Value = "Etiopia"
$ws.Range("B141").Value = 124
$ws.Range("C141").Value = 1
$ws.Range("D141").Value = 50
$ws.Range("E141").Value = 71
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 3

# the old Etiopia row has been pushed down to row 143 - remove it, leaving
# Madagascar (now row 142, figures unchanged) directly below the new Etiopia row
$ws.Rows("143").Delete()
